$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.681921601295471
$ws.Range("B1").Value = 3.644757747650146
$ws.Range("C1").Value = 3.20806884765625
$ws.Range("D1").Value = 3.476157903671265
$ws.Range("E1").Value = 1.616185307502747
